# Updated cryptos list on Wed Oct 30 17:54:05 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a value as plain text, avoiding Excel's automatic
# number/date inference (e.g. "598.01" or "71.89" would otherwise be
# silently converted into a numeric value). We temporarily force a text
# number format while assigning the value, then restore the cell's style
# back to the workbook's default "Normal" style so no stray formatting
# is introduced into the saved file.
function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "71.833.98"
Set-TextCell "E2" "  -1.06%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.673.49"
Set-TextCell "E3" "  +0.60%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  -0.02%  "

# Row 5 - BNB
Set-TextCell "D5" "598.01"
Set-TextCell "E5" "  -1.76%  "

# Row 6 - Solana
Set-TextCell "D6" "175.08"
Set-TextCell "E6" "  -3.00%  "

# Row 7 - USDC
Set-TextCell "E7" "  -0.06%  "

# Row 8 - XRP
Set-TextCell "D8" "0.523"
Set-TextCell "E8" "  -0.97%  "

# Row 9 - LidoStakedEther
Set-TextCell "D9" "2.672.84"
Set-TextCell "E9" "  +0.58%  "

# Row 10 - Dogecoin
Set-TextCell "E10" "  -5.10%  "

# Row 11 - TRON
Set-TextCell "E11" "  +2.06%  "

# Row 12 - Cardano
Set-TextCell "D12" "0.356"
Set-TextCell "E12" "  +0.77%  "

# Row 13 - Toncoin
Set-TextCell "E13" "  -2.08%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextCell "D14" "3.170.54"

# Row 15 - ShibaInu
Set-TextCell "D15" "0.0000185"
Set-TextCell "E15" "  -4.34%  "

# Row 16 - WrappedBTC
Set-TextCell "D16" "71.757.52"
Set-TextCell "E16" "  -1.10%  "

# Row 17 - Avalanche
Set-TextCell "D17" "26.19"
Set-TextCell "E17" "  -3.00%  "

# Row 18 - WrappedEther
Set-TextCell "D18" "2.661.61"
Set-TextCell "E18" "  +0.11%  "

# Row 19 - Chainlink
Set-TextCell "D19" "12.21"
Set-TextCell "E19" "  +5.36%  "

# Row 20 - Uniswap
Set-TextCell "D20" "8.26"
Set-TextCell "E20" "  +4.12%  "

# Row 21 - BitcoinCash
Set-TextCell "D21" "370.85"
Set-TextCell "E21" "  -3.46%  "

# Row 22 - Polkadot
Set-TextCell "D22" "4.16"
Set-TextCell "E22" "  -1.41%  "

# Row 23 - SuiNetwork
Set-TextCell "E23" "  -1.47%  "

# Row 24 - Litecoin
Set-TextCell "D24" "72.00"
Set-TextCell "E24" "  -1.49%  "

# Row 25 - Dai
Set-TextCell "E25" "  +0.05%  "

# Row 26 - NEARProtocol
Set-TextCell "D26" "4.33"
Set-TextCell "E26" "  -2.71%  "

# Row 27 - Aptos
Set-TextCell "D27" "9.76"
Set-TextCell "E27" "  -2.52%  "

# Row 28 - WrappedeETH
Set-TextCell "D28" "2.816.69"
Set-TextCell "E28" "  +0.67%  "

# Row 29 - Binance-PegBSC-USD
Set-TextCell "E29" "  -0.07%  "

# Row 30 - PEPE
Set-TextCell "D30" "0.0₃0970"
Set-TextCell "E30" "  -0.58%  "

# Row 31 - InternetComputer(DFINITY)
Set-TextCell "D31" "8.05"
Set-TextCell "E31" "  -0.70%  "

# Row 32 - Bittensor
Set-TextCell "D32" "501.83"
Set-TextCell "E32" "  -7.98%  "

# Row 33 - Fetch.AI
Set-TextCell "D33" "1.30"
Set-TextCell "E33" "  -3.72%  "

# Row 34 - PancakeSwap
Set-TextCell "E34" "  -1.72%  "

# Row 35 - FirstDigitalUSD
Set-TextCell "E35" "  +0.00%  "

# Row 36 - Monero
Set-TextCell "D36" "162.27"
Set-TextCell "E36" "  -1.94%  "

# Row 37 - EthereumClassic
Set-TextCell "D37" "19.55"
Set-TextCell "E37" "  +0.78%  "

# Row 38 - WhiteBITCoin
Set-TextCell "D38" "19.07"
Set-TextCell "E38" "  -0.31%  "

# Row 39 - ImmutableX
Set-TextCell "D39" "1.38"
Set-TextCell "E39" "  -2.98%  "

# Row 40 - Kaspa
Set-TextCell "D40" "0.110"
Set-TextCell "E40" "  -3.96%  "

# Row 41 - Stacks
Set-TextCell "E41" "  -5.01%  "

# Row 42 - USDe
Set-TextCell "E42" "  -0.15%  "

# Row 43 - RenderToken
Set-TextCell "E43" "  -2.41%  "

# Row 44 - dogwifhat
Set-TextCell "D44" "2.57"
Set-TextCell "E44" "  -2.91%  "

# Row 45 - PolygonEcosystemToken
Set-TextCell "E45" "  -1.13%  "

# Row 46 - was OKB, now Aave
Set-TextCell "B46" "Aave"
Set-TextCell "C46" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell "D46" "156.06"
Set-TextCell "E46" "  +1.93%  "

# Row 47 - was Aave, now OKB
Set-TextCell "B47" "OKB"
Set-TextCell "C47" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D47" "39.48"
Set-TextCell "E47" "  -0.79%  "

# Row 48 - ARBITRUM
Set-TextCell "D48" "0.562"
Set-TextCell "E48" "  +2.96%  "

# Row 49 - Filecoin
Set-TextCell "E49" "  +0.12%  "

# Row 50 - Optimism
Set-TextCell "E50" "  +1.32%  "

# Row 51 - Cronos
Set-TextCell "E51" "  -1.51%  "
